$d = $word.ActiveDocument

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# ---------------------------------------------------------------------------
# 1) Title-block date line: split "Effective: October 31, 2025 | v1.0" so the
#    year "2025" becomes its own run wrapped in gramStart/gramEnd proofErr
#    markers (mirrors Word's grammar-check run splitting).
# ---------------------------------------------------------------------------
$range1 = $d.Content
$found1 = $range1.Find.Execute("Effective: October 31, 2025 | v1.0")
if ($found1) {
    $r1 = $d.Range($range1.Start, $range1.End)
    $xml1 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="0D3F5286" w14:textId="77777777" w:rsidR="00D678B1" w:rsidRDefault="00000000"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Effective: October 31, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>2025</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> | v1.0</w:t></w:r></w:p>'
    $r1.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# 2) L1 - Associate table cell: split "Learns core skills" into "Learns" /
#    " core skills" wrapped in gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$range2 = $d.Content
$found2 = $range2.Find.Execute("Learns core skills")
if ($found2) {
    $r2 = $d.Range($range2.Start, $range2.End)
    $xml2 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="27A6C6FA" w14:textId="77777777" w:rsidR="00D678B1" w:rsidRDefault="00000000"><w:proofErr w:type="gramStart"/><w:r><w:t>Learns</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> core skills</w:t></w:r></w:p>'
    $r2.InsertXML($xml2)
}

# ---------------------------------------------------------------------------
# 3) FAQ bullet: split "Q: Can I skip a level? A: Rarely; requires sustained
#    evidence of operating two levels up." so "Rarely;" becomes its own run
#    wrapped in gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$range3 = $d.Content
$found3 = $range3.Find.Execute("Q: Can I skip a level? A: Rarely; requires sustained evidence of operating two levels up.")
if ($found3) {
    $r3 = $d.Range($range3.Start, $range3.End)
    $xml3 = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="3E3AF761" w14:textId="77777777" w:rsidR="00D678B1" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t xml:space="preserve">Q: Can I skip a level? A: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Rarely;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> requires sustained evidence of operating two levels up.</w:t></w:r></w:p>'
    $r3.InsertXML($xml3)
}

Write-Host "done: found1=$found1 found2=$found2 found3=$found3"
